$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: add Q8 label in J1 (match style/format of existing header cells)
$ws.Range("J1").Value = "Q8"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = "2022-01-01 00:00:00_diff"
$ws.Range("B2").Value = -0.01758111206922311
$ws.Range("C2").Value = -0.9703086835434362
$ws.Range("D2").Value = -1.190204369659767
$ws.Range("E2").Value = 0.851657034662153
$ws.Range("F2").Value = 1.077585438625931
$ws.Range("G2").Value = -0.1996304584251192
$ws.Range("H2").Value = 0.1134005067055028

# Row 3
$ws.Range("A3").Value = "2022-04-01 00:00:00_diff"
$ws.Range("B3").Value = -0.8029421940374804
$ws.Range("C3").Value = -1.022837880153811
$ws.Range("D3").Value = 1.019023524168109
$ws.Range("E3").Value = 1.244951928131887
$ws.Range("F3").Value = -0.0322639689191633
$ws.Range("G3").Value = 0.2807669962114587

# Row 4
$ws.Range("A4").Value = "2022-07-01 00:00:00_diff"
$ws.Range("B4").Value = -0.7936197797752114
$ws.Range("C4").Value = 1.248241624546709
$ws.Range("D4").Value = 1.474170028510487
$ws.Range("E4").Value = 0.1969541314594366
$ws.Range("F4").Value = 0.5099850965900585
$ws.Range("G4").Value = 0.007569982732279296
$ws.Range("H4").Value = 0.6592032301411037
$ws.Range("I4").Value = 0.4663036110991228
$ws.Range("J4").Value = -0.04443783748577212

# Row 5
$ws.Range("A5").Value = "2022-10-01 00:00:00_diff"
$ws.Range("B5").Value = 2.484163736993811
$ws.Range("C5").Value = 2.710092140957589
$ws.Range("D5").Value = 1.432876243906539
$ws.Range("E5").Value = 1.745907209037161
$ws.Range("F5").Value = 1.243492095179382
$ws.Range("G5").Value = 1.895125342588206
$ws.Range("H5").Value = 1.702225723546225
$ws.Range("I5").Value = 1.19148427496133

# Row 6
$ws.Range("A6").Value = "2023-01-01 00:00:00_diff"
$ws.Range("B6").Value = 1.371380565536508
$ws.Range("C6").Value = 0.09416466848545757
$ws.Range("D6").Value = 0.4071956336160796
$ws.Range("E6").Value = -0.09521948024169971
$ws.Range("F6").Value = 0.5564137671671248
$ws.Range("G6").Value = 0.3635141481251438
$ws.Range("H6").Value = -0.1472273004597511

# Row 7
$ws.Range("A7").Value = "2023-04-01 00:00:00_diff"
$ws.Range("B7").Value = 0.2659007569564139
$ws.Range("C7").Value = 0.5789317220870359
$ws.Range("D7").Value = 0.07651660822925663
$ws.Range("E7").Value = 0.7281498556380811
$ws.Range("F7").Value = 0.5352502365961002
$ws.Range("G7").Value = 0.02450878801120521

# Row 8
$ws.Range("A8").Value = "2023-07-01 00:00:00_diff"
$ws.Range("B8").Value = 1.260690851164143
$ws.Range("C8").Value = 0.7582757373063643
$ws.Range("D8").Value = 1.409908984715189
$ws.Range("E8").Value = 1.217009365673208
$ws.Range("F8").Value = 0.7062679170883128
$ws.Range("G8").Value = 0.911668649685511
$ws.Range("H8").Value = 0.5051827077222001
$ws.Range("I8").Value = 0.9402868649905415

# Row 9
$ws.Range("A9").Value = "2023-10-01 00:00:00_diff"
$ws.Range("B9").Value = 0.2946970959196917
$ws.Range("C9").Value = 0.9463303433285162
$ws.Range("D9").Value = 0.7534307242865352
$ws.Range("E9").Value = 0.2426892757016403
$ws.Range("F9").Value = 0.4480900082988384
$ws.Range("G9").Value = 0.04160406633552749
$ws.Range("H9").Value = 0.4767082236038689

# Row 10
$ws.Range("A10").Value = "2024-01-01 00:00:00_diff"
$ws.Range("B10").Value = 0.3856725119803543
$ws.Range("C10").Value = 0.1927728929383733
$ws.Range("D10").Value = -0.3179685556465216
$ws.Range("E10").Value = -0.1125678230493235
$ws.Range("F10").Value = -0.5190537650126344
$ws.Range("G10").Value = -0.08394960774429301

# Row 11
$ws.Range("A11").Value = "2024-04-01 00:00:00_diff"
$ws.Range("B11").Value = 0.2952882579329085
$ws.Range("C11").Value = -0.2154531906519864
$ws.Range("D11").Value = -0.01005245805478834
$ws.Range("E11").Value = -0.4165384000180992
$ws.Range("F11").Value = 0.01856575725024216

# Row 12
$ws.Range("A12").Value = "2024-07-01 00:00:00_diff"
$ws.Range("B12").Value = -0.2970557949068323
$ws.Range("C12").Value = -0.09165506230963413
$ws.Range("D12").Value = -0.4981410042729451
$ws.Range("E12").Value = -0.06303684700460363

# Row 13
$ws.Range("A13").Value = "2024-10-01 00:00:00_diff"
$ws.Range("B13").Value = -0.1069508448768545
$ws.Range("C13").Value = -0.5134367868401654
$ws.Range("D13").Value = -0.07833262957182399

# Row 14
$ws.Range("A14").Value = "2025-01-01 00:00:00_diff"
$ws.Range("B14").Value = -0.6053253388254292
$ws.Range("C14").Value = -0.1702211815570877

# Row 15
$ws.Range("A15").Value = "2025-04-01 00:00:00_diff"
$ws.Range("B15").Value = 0.06843616378760228

# Row 16
$ws.Range("A16").Value = "2025-07-01 00:00:00_diff"
